$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "CreatedAt: 2025-08-29T19:07:13"
$ws.Range("W4").Value = 72.34999999999999
$ws.Range("X4").Value = 47.2
$ws.Range("Y4").Value = -2.54
$ws.Range("Z4").Value = 34.94
$ws.Range("X6").Value = -1.09
$ws.Range("Y6").Value = -1.01
$ws.Range("Z6").Value = -0.87
$ws.Range("X8").Value = 10.86
$ws.Range("Y8").Value = -37.51
$ws.Range("X9").Value = 48.25
$ws.Range("Y9").Value = -1.74
$ws.Range("Z9").Value = 35.78
$ws.Range("W11").Value = -1.18
$ws.Range("X11").Value = -0.04
$ws.Range("Y11").Value = -0.21
$ws.Range("Z11").Value = -0.04
$ws.Range("W12").Value = 26.01
$ws.Range("X13").Value = 10.86
$ws.Range("Y13").Value = -37.51
$ws.Range("X14").Value = 48.25
$ws.Range("Y14").Value = -1.7
$ws.Range("Z14").Value = 35.78
$ws.Range("W16").Value = -1.18
$ws.Range("X16").Value = -0.04
$ws.Range("Y16").Value = -0.18
$ws.Range("Z16").Value = -0.04
$ws.Range("W17").Value = 26.01
$ws.Range("X18").Value = 10.86
$ws.Range("Y18").Value = -37.51
$ws.Range("W19").Value = 72.22
$ws.Range("X19").Value = 47.13
$ws.Range("Y19").Value = -2.57
$ws.Range("Z19").Value = 34.91
$ws.Range("X21").Value = -1.16
$ws.Range("Y21").Value = -1.05
$ws.Range("Z21").Value = -0.91
$ws.Range("X23").Value = 10.86
$ws.Range("Y23").Value = -37.51
$ws.Range("W24").Value = 72.22
$ws.Range("X24").Value = 47.13
$ws.Range("Y24").Value = -2.57
$ws.Range("Z24").Value = 34.91
$ws.Range("X26").Value = -1.16
$ws.Range("Y26").Value = -1.05
$ws.Range("Z26").Value = -0.91
$ws.Range("X28").Value = 10.86
$ws.Range("Y28").Value = -37.51
$ws.Range("W29").Value = 72.08
$ws.Range("X29").Value = 47.02
$ws.Range("Y29").Value = -2.6
$ws.Range("Z29").Value = 34.88
$ws.Range("W31").Value = -3.1
$ws.Range("X31").Value = -1.27
$ws.Range("Y31").Value = -1.08
$ws.Range("Z31").Value = -0.9399999999999999
$ws.Range("X33").Value = 10.86
$ws.Range("Y33").Value = -37.51
$ws.Range("W34").Value = 74.43000000000001
$ws.Range("X34").Value = 48.85
$ws.Range("Y34").Value = -1.27
$ws.Range("Z34").Value = 36.25
$ws.Range("X36").Value = 0.53
$ws.Range("Y36").Value = 0.25
$ws.Range("Z36").Value = 0.44
$ws.Range("X37").Value = 0.02
$ws.Range("X38").Value = 10.86
$ws.Range("Y38").Value = -37.51
$ws.Range("W39").Value = 72.34999999999999
$ws.Range("X39").Value = 47.2
$ws.Range("Y39").Value = -2.54
$ws.Range("Z39").Value = 34.94
$ws.Range("X41").Value = -1.09
$ws.Range("Y41").Value = -1.01
$ws.Range("Z41").Value = -0.87
$ws.Range("X43").Value = 10.86
$ws.Range("Y43").Value = -37.51
$ws.Range("W44").Value = 76.03
$ws.Range("X44").Value = 48.73
$ws.Range("Y44").Value = -1.17
$ws.Range("Z44").Value = 36.39
$ws.Range("W46").Value = 0.85
$ws.Range("X46").Value = 0.44
$ws.Range("Y46").Value = 0.35
$ws.Range("Z46").Value = 0.57
$ws.Range("X48").Value = 10.86
$ws.Range("Y48").Value = -37.51
$ws.Range("W49").Value = 77.5
$ws.Range("X49").Value = 49.05
$ws.Range("Y49").Value = -1.01
$ws.Range("Z49").Value = 36.55
$ws.Range("W51").Value = 2.33
$ws.Range("Y51").Value = 0.51
$ws.Range("Z51").Value = 0.73
$ws.Range("X53").Value = 10.86
$ws.Range("Y53").Value = -37.51
$ws.Range("W54").Value = 76.40000000000001
$ws.Range("X54").Value = 49.01
$ws.Range("Y54").Value = -0.64
$ws.Range("Z54").Value = 36.7
$ws.Range("W56").Value = 1.22
$ws.Range("X56").Value = 0.73
$ws.Range("Y56").Value = 0.89
$ws.Range("Z56").Value = 0.88
$ws.Range("X58").Value = 10.86
$ws.Range("Y58").Value = -37.51
$ws.Range("W59").Value = 78.88
$ws.Range("X59").Value = 50.09
$ws.Range("Y59").Value = 0.14
$ws.Range("Z59").Value = 37.62
$ws.Range("W61").Value = 3.71
$ws.Range("X61").Value = 1.8
$ws.Range("Y61").Value = 1.66
$ws.Range("Z61").Value = 1.81
$ws.Range("X63").Value = 10.86
$ws.Range("Y63").Value = -37.51
$ws.Range("W64").Value = 79.72
$ws.Range("X64").Value = 50.51
$ws.Range("Y64").Value = 0.49
$ws.Range("Z64").Value = 37.9
$ws.Range("W66").Value = 4.54
$ws.Range("X66").Value = 2.22
$ws.Range("Y66").Value = 2.01
$ws.Range("Z66").Value = 2.08
$ws.Range("X68").Value = 10.86
$ws.Range("Y68").Value = -37.51
$ws.Range("W69").Value = 80.66
$ws.Range("X69").Value = 50.98
$ws.Range("Y69").Value = 0.98
$ws.Range("Z69").Value = 38.47
$ws.Range("W71").Value = 5.49
$ws.Range("X71").Value = 2.69
$ws.Range("Y71").Value = 2.5
$ws.Range("Z71").Value = 2.65
$ws.Range("X73").Value = 10.86
$ws.Range("Y73").Value = -37.51
$ws.Range("W74").Value = 78.15000000000001
$ws.Range("X74").Value = 49.77
$ws.Range("Y74").Value = -0.22
$ws.Range("Z74").Value = 37.27
$ws.Range("W76").Value = 2.97
$ws.Range("X76").Value = 1.48
$ws.Range("Y76").Value = 1.31
$ws.Range("Z76").Value = 1.45
$ws.Range("X78").Value = 10.86
$ws.Range("Y78").Value = -37.51
$ws.Range("W79").Value = 78.61
$ws.Range("X79").Value = 50.01
$ws.Range("Y79").Value = 0.02
$ws.Range("Z79").Value = 37.5
$ws.Range("W81").Value = 3.44
$ws.Range("X81").Value = 1.72
$ws.Range("Y81").Value = 1.54
$ws.Range("Z81").Value = 1.68
$ws.Range("X83").Value = 10.86
$ws.Range("Y83").Value = -37.51
$ws.Range("W84").Value = 76.48
$ws.Range("X84").Value = 49.01
$ws.Range("Y84").Value = -0.6
$ws.Range("Z84").Value = 36.81
$ws.Range("W86").Value = 1.3
$ws.Range("X86").Value = 0.73
$ws.Range("Y86").Value = 0.92
$ws.Range("Z86").Value = 0.99
$ws.Range("X88").Value = 10.86
$ws.Range("Y88").Value = -37.51
$ws.Range("W89").Value = 72.08
$ws.Range("X89").Value = 47.02
$ws.Range("Y89").Value = -2.6
$ws.Range("Z89").Value = 34.88
$ws.Range("W91").Value = -3.1
$ws.Range("X91").Value = -1.27
$ws.Range("Y91").Value = -1.08
$ws.Range("Z91").Value = -0.9399999999999999
$ws.Range("X93").Value = 10.86
$ws.Range("Y93").Value = -37.51
